$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 90.56466633333334
$ws.Range("H2").Value = 271.693999
$ws.Range("I2").Value = 0.2234788625831797
$ws.Range("J2").Value = 0.2234788625831796
$ws.Range("M2").Value = 3.456265333333333
$ws.Range("N2").Value = 10.368796
$ws.Range("O2").Value = 0.009841535807677501
$ws.Range("P2").Value = 0.0098415358076775
$ws.Range("Q2").Value = 313.0155166728004
$ws.Range("R2").Value = 2817.139650055204
$ws.Range("S2").Value = 0.002199375228371402
$ws.Range("T2").Value = 0.002199375228371402
$ws.Range("G3").Value = 90.56466633333334
$ws.Range("H3").Value = 271.693999
$ws.Range("I3").Value = 0.2234788625831797
$ws.Range("J3").Value = 0.2234788625831796
$ws.Range("O3").Value = 0.8587907398420774
$ws.Range("P3").Value = 0.8587907398420773
$ws.Range("Q3").Value = 27314.3168300804
$ws.Range("R3").Value = 245828.8514707236
$ws.Range("S3").Value = 0.1919215777368748
$ws.Range("T3").Value = 0.1919215777368748
$ws.Range("G4").Value = 90.56466633333334
$ws.Range("H4").Value = 271.693999
$ws.Range("I4").Value = 0.2234788625831797
$ws.Range("J4").Value = 0.2234788625831796
$ws.Range("O4").Value = 0.1313677243502452
$ws.Range("P4").Value = 0.1313677243502452
$ws.Range("Q4").Value = 4178.223492266695
$ws.Range("R4").Value = 37604.01143040026
$ws.Range("S4").Value = 0.02935790961793347
$ws.Range("T4").Value = 0.02935790961793346
$ws.Range("I5").Value = 0.601197186834308
$ws.Range("J5").Value = 0.6011971868343079
$ws.Range("M5").Value = 3.456265333333333
$ws.Range("N5").Value = 10.368796
$ws.Range("O5").Value = 0.009841535807677501
$ws.Range("P5").Value = 0.0098415358076775
$ws.Range("Q5").Value = 842.066430283232
$ws.Range("R5").Value = 7578.597872549089
$ws.Range("S5").Value = 0.005916703641704824
$ws.Range("T5").Value = 0.005916703641704821
$ws.Range("I6").Value = 0.601197186834308
$ws.Range("J6").Value = 0.6011971868343079
$ws.Range("O6").Value = 0.8587907398420774
$ws.Range("P6").Value = 0.8587907398420773
$ws.Range("S6").Value = 0.516302576872411
$ws.Range("T6").Value = 0.5163025768724109
$ws.Range("I7").Value = 0.601197186834308
$ws.Range("J7").Value = 0.6011971868343079
$ws.Range("O7").Value = 0.1313677243502452
$ws.Range("P7").Value = 0.1313677243502452
$ws.Range("S7").Value = 0.07897790632019223
$ws.Range("T7").Value = 0.07897790632019219
$ws.Range("I8").Value = 0.1753239505825123
$ws.Range("J8").Value = 0.1753239505825123
$ws.Range("M8").Value = 3.456265333333333
$ws.Range("N8").Value = 10.368796
$ws.Range("O8").Value = 0.009841535807677501
$ws.Range("P8").Value = 0.0098415358076775
$ws.Range("Q8").Value = 245.5673719758413
$ws.Range("R8").Value = 2210.106347782572
$ws.Range("S8").Value = 0.001725456937601275
$ws.Range("T8").Value = 0.001725456937601275
$ws.Range("I9").Value = 0.1753239505825123
$ws.Range("J9").Value = 0.1753239505825123
$ws.Range("O9").Value = 0.8587907398420774
$ws.Range("P9").Value = 0.8587907398420773
$ws.Range("S9").Value = 0.1505665852327916
$ws.Range("T9").Value = 0.1505665852327915
$ws.Range("I10").Value = 0.1753239505825123
$ws.Range("J10").Value = 0.1753239505825123
$ws.Range("O10").Value = 0.1313677243502452
$ws.Range("P10").Value = 0.1313677243502452
$ws.Range("S10").Value = 0.02303190841211948
$ws.Range("T10").Value = 0.02303190841211947
